$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-31 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-01 Saturday", 2)

$d.Content.Find.Execute("788×3=2364", $true, $false, $false, $false, $false, $true, 1, $false, "335×9=3015", 2)
$d.Content.Find.Execute("308×5=1540", $true, $false, $false, $false, $false, $true, 1, $false, "131×5=655", 2)
$d.Content.Find.Execute("634×5=3170", $true, $false, $false, $false, $false, $true, 1, $false, "895×2=1790", 2)
$d.Content.Find.Execute("170×3=510", $true, $false, $false, $false, $false, $true, 1, $false, "420×2=840", 2)
$d.Content.Find.Execute("137×6=822", $true, $false, $false, $false, $false, $true, 1, $false, "205×8=1640", 2)

$d.Content.Find.Execute("831×2=1662", $true, $false, $false, $false, $false, $true, 1, $false, "686×3=2058", 2)
$d.Content.Find.Execute("122×3=366", $true, $false, $false, $false, $false, $true, 1, $false, "704×9=6336", 2)
$d.Content.Find.Execute("304×3=912", $true, $false, $false, $false, $false, $true, 1, $false, "403×2=806", 2)
$d.Content.Find.Execute("581×4=2324", $true, $false, $false, $false, $false, $true, 1, $false, "917×4=3668", 2)
$d.Content.Find.Execute("991×6=5946", $true, $false, $false, $false, $false, $true, 1, $false, "654×9=5886", 2)

$d.Content.Find.Execute("928×2=1856", $true, $false, $false, $false, $false, $true, 1, $false, "570×4=2280", 2)
$d.Content.Find.Execute("346×7=2422", $true, $false, $false, $false, $false, $true, 1, $false, "840×5=4200", 2)
$d.Content.Find.Execute("402×5=2010", $true, $false, $false, $false, $false, $true, 1, $false, "189×3=567", 2)
$d.Content.Find.Execute("259×3=777", $true, $false, $false, $false, $false, $true, 1, $false, "588×8=4704", 2)
$d.Content.Find.Execute("378×7=2646", $true, $false, $false, $false, $false, $true, 1, $false, "262×6=1572", 2)

$d.Content.Find.Execute("840×3=2520", $true, $false, $false, $false, $false, $true, 1, $false, "181×6=1086", 2)
$d.Content.Find.Execute("205×5=1025", $true, $false, $false, $false, $false, $true, 1, $false, "578×9=5202", 2)
$d.Content.Find.Execute("230×7=1610", $true, $false, $false, $false, $false, $true, 1, $false, "669×4=2676", 2)
$d.Content.Find.Execute("299×7=2093", $true, $false, $false, $false, $false, $true, 1, $false, "790×3=2370", 2)
$d.Content.Find.Execute("380×4=1520", $true, $false, $false, $false, $false, $true, 1, $false, "411×9=3699", 2)

$d.Content.Find.Execute("203×4=812", $true, $false, $false, $false, $false, $true, 1, $false, "323×9=2907", 2)
$d.Content.Find.Execute("968×4=3872", $true, $false, $false, $false, $false, $true, 1, $false, "818×7=5726", 2)
$d.Content.Find.Execute("798×9=7182", $true, $false, $false, $false, $false, $true, 1, $false, "792×6=4752", 2)
$d.Content.Find.Execute("901×3=2703", $true, $false, $false, $false, $false, $true, 1, $false, "820×3=2460", 2)
$d.Content.Find.Execute("958×8=7664", $true, $false, $false, $false, $false, $true, 1, $false, "686×9=6174", 2)
